$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-5 with new Node Identifier / SH / SL values
$ws.Range("A2").Value = "DM-2"
$ws.Range("B2").Value = "0013A200"
$ws.Range("C2").Value = "40F96362"

$ws.Range("A3").Value = "DMPCB-3"
$ws.Range("B3").Value = "0013A200"
$ws.Range("C3").Value = "40A863C2"

$ws.Range("A4").Value = "DMPCB-4"
$ws.Range("B4").Value = "0013A200"
$ws.Range("C4").Value = "40A164C9"

$ws.Range("A5").Value = "DM-5"
$ws.Range("B5").Value = "0013A200"
$ws.Range("C5").Value = "40F96363"

# Add new row 6
$ws.Range("A6").Value = "DM-6"
$ws.Range("B6").Value = "0013A200"
$ws.Range("C6").Value = "40F96376"

# Update selection to match the diff (C9)
$ws.Range("C9").Select()
